# Activities Test data changes - 13 Dec 2023
# Update the "Users" sheet: rename the user in A2 from "Drew Koecher" to
# "Ayati Arvind", and leave the sheet's selection on A3 (matching the
# cursor position Excel persisted when the file was last saved).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")
$ws.Activate()

$ws.Range("A2").Value = "Ayati Arvind"

$ws.Range("A3").Select()
